# Apply the "Updated symbol list" edit to the cryptos worksheet.
# All target cells store values as text (inline strings) in the original
# workbook, so we force the cell format to Text ("@") before writing the
# new values to avoid Excel auto-converting numeric-looking strings into
# real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "245.84"
Set-TextValue "D3"  "21.97"
Set-TextValue "D4"  "5.409"
Set-TextValue "D5"  "0.05761"
Set-TextValue "D6"  "3.401"
Set-TextValue "D7"  "6.333"
Set-TextValue "D8"  "0.8145"
Set-TextValue "D9"  "0.9859"
Set-TextValue "D10" "0.1434"
Set-TextValue "D11" "0.07411"
Set-TextValue "D12" "0.03140"
Set-TextValue "D13" "0.03017"
Set-TextValue "D14" "4.138"
Set-TextValue "D15" "0.09402"
Set-TextValue "D16" "0.001589"
Set-TextValue "D17" "0.04819"
Set-TextValue "D18" "0.0005848"
Set-TextValue "D19" "0.006217"
Set-TextValue "D20" "0.004111"
Set-TextValue "D21" "0.0009948"
Set-TextValue "D23" "3.759"
Set-TextValue "D24" "2.204"
Set-TextValue "D25" "0.3257"
Set-TextValue "D26" "0.1329"
Set-TextValue "D40" "0.03887"

# --- Row 41: KickToken -> BKEXToken ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1074"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# --- Row 42: BKEXToken -> CEJI ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002655"
$ws.Range("E42").Value = "41CEJICEJI"

# --- Row 43: CEJI -> KickToken ---
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006448"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Remaining price (column D) updates ---
Set-TextValue "D44" "0.006279"
Set-TextValue "D45" "0.00005593"
Set-TextValue "D49" "0.00002099"
